$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 21.619047
$ws.Range("I11").Value = 21.619047
$ws.Range("K11").Value = 21.619047
$ws.Range("M11").Value = 118.380953

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1495.1428
$ws.Range("J19").Value = 1693.2
$ws.Range("L19").Value = 1693.2
$ws.Range("N19").Value = -2043.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1877.4445
$ws.Range("J40").Value = 1877.4445
$ws.Range("L40").Value = 1877.4445
$ws.Range("N40").Value = -2227.4445

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 17701.5
$ws.Range("J43").Value = 1055
$ws.Range("L43").Value = 1055
$ws.Range("N43").Value = -1193

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1321.5834
$ws.Range("I98").Value = 1215.4286
$ws.Range("J98").Value = 2064.6667
$ws.Range("K98").Value = 1215.4286
$ws.Range("L98").Value = 2064.6667
$ws.Range("M98").Value = 282.5714
$ws.Range("N98").Value = -5060.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1094.2858
$ws.Range("J111").Value = 930
$ws.Range("L111").Value = 2790
$ws.Range("N111").Value = -8924

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 56432.91
$ws.Range("I113").Value = 110280
$ws.Range("J113").Value = 11560.333
$ws.Range("K113").Value = 110280
$ws.Range("L113").Value = 11560.333
$ws.Range("M113").Value = -107026
$ws.Range("N113").Value = -18068.333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1321.5834
$ws.Range("I122").Value = 1215.4286
$ws.Range("J122").Value = 2064.6667
$ws.Range("K122").Value = 3646.2858
$ws.Range("L122").Value = 6194.000100000001
$ws.Range("M122").Value = -1196.2858
$ws.Range("N122").Value = -11094.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1677.1428
$ws.Range("I135").Value = 1290
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 11610
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -9075
$ws.Range("N135").Value = -41070

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1399.2632
$ws.Range("I137").Value = 1417.9375
$ws.Range("K137").Value = 4253.8125
$ws.Range("M137").Value = -1703.8125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3150.1975
$ws.Range("I138").Value = 1857.9706
$ws.Range("J138").Value = 4085
$ws.Range("K138").Value = 5573.9118
$ws.Range("L138").Value = 12255
$ws.Range("M138").Value = -433.9117999999999
$ws.Range("N138").Value = -22535

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1907.4375
$ws.Range("I141").Value = 1679.9286
$ws.Range("J141").Value = 3500
$ws.Range("K141").Value = 5039.7858
$ws.Range("L141").Value = 10500
$ws.Range("M141").Value = 140.2142000000003
$ws.Range("N141").Value = -20860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9798.5625
$ws.Range("I45").Value = 9987.200000000001
$ws.Range("K45").Value = 9987.200000000001
$ws.Range("M45").Value = -9610.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 500224.5
$ws.Range("J69").Value = 500224.5
$ws.Range("L69").Value = 500224.5
$ws.Range("N69").Value = -501722.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H72").Value = 500224.5
$ws.Range("J72").Value = 500224.5
$ws.Range("L72").Value = 1500673.5
$ws.Range("N72").Value = -1508161.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6399.722
$ws.Range("I74").Value = 1376.4333
$ws.Range("J74").Value = 31516.166
$ws.Range("K74").Value = 1376.4333
$ws.Range("L74").Value = 31516.166
$ws.Range("M74").Value = -502.4332999999999
$ws.Range("N74").Value = -33264.166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6399.722
$ws.Range("I77").Value = 1376.4333
$ws.Range("J77").Value = 31516.166
$ws.Range("K77").Value = 6882.166499999999
$ws.Range("L77").Value = 157580.83
$ws.Range("M77").Value = -2514.166499999999
$ws.Range("N77").Value = -166316.83

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1025.4828
$ws.Range("I97").Value = 783.1667
$ws.Range("K97").Value = 783.1667
$ws.Range("M97").Value = -287.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 80721
$ws.Range("I133").Value = 33000
$ws.Range("J133").Value = 85493.10000000001
$ws.Range("K133").Value = 33000
$ws.Range("L133").Value = 85493.10000000001
$ws.Range("M133").Value = -30470
$ws.Range("N133").Value = -90553.10000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 70000
$ws.Range("J27").Value = 70000
$ws.Range("L27").Value = 70000
$ws.Range("N27").Value = -70384

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3893.6191
$ws.Range("I99").Value = 1149.1333
$ws.Range("K99").Value = 1149.1333
$ws.Range("M99").Value = 348.8667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2206.366
$ws.Range("I134").Value = 1783.6154
$ws.Range("K134").Value = 5350.8462
$ws.Range("M134").Value = -2815.8462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 98
$ws.Range("I2").Value = 98
$ws.Range("K2").Value = 98
$ws.Range("M2").Value = 15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4563.9165
$ws.Range("I31").Value = 3599.5
$ws.Range("J31").Value = 5046.125
$ws.Range("K31").Value = 3599.5
$ws.Range("L31").Value = 5046.125
$ws.Range("M31").Value = -3304.5
$ws.Range("N31").Value = -5636.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4563.9165
$ws.Range("I34").Value = 3599.5
$ws.Range("J34").Value = 5046.125
$ws.Range("K34").Value = 3599.5
$ws.Range("L34").Value = 5046.125
$ws.Range("M34").Value = -3397.5
$ws.Range("N34").Value = -5450.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2687.3953
$ws.Range("I58").Value = 2917.3572
$ws.Range("J58").Value = 2258.1333
$ws.Range("K58").Value = 2917.3572
$ws.Range("L58").Value = 2258.1333
$ws.Range("M58").Value = -2714.3572
$ws.Range("N58").Value = -2664.1333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5299.2
$ws.Range("I99").Value = 2570.2856
$ws.Range("K99").Value = 2570.2856
$ws.Range("M99").Value = -1072.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5299.2
$ws.Range("I126").Value = 2570.2856
$ws.Range("K126").Value = 7710.8568
$ws.Range("M126").Value = -5240.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4122.1724
$ws.Range("I132").Value = 4111.0435
$ws.Range("J132").Value = 4164.8335
$ws.Range("K132").Value = 12333.1305
$ws.Range("L132").Value = 12494.5005
$ws.Range("M132").Value = -9803.130499999999
$ws.Range("N132").Value = -17554.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 11731.974
$ws.Range("I134").Value = 5177.5625
$ws.Range("J134").Value = 46688.832
$ws.Range("K134").Value = 15532.6875
$ws.Range("L134").Value = 140066.496
$ws.Range("M134").Value = -12997.6875
$ws.Range("N134").Value = -145136.496

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2687.3953
$ws.Range("I136").Value = 2917.3572
$ws.Range("J136").Value = 2258.1333
$ws.Range("K136").Value = 8752.071599999999
$ws.Range("L136").Value = 6774.3999
$ws.Range("M136").Value = -6202.071599999999
$ws.Range("N136").Value = -11874.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 59611744
$ws.Range("I4").Value = 70695170
$ws.Range("J4").Value = 500133
$ws.Range("K4").Value = 212085510
$ws.Range("L4").Value = 1500399
$ws.Range("M4").Value = -212085398
$ws.Range("N4").Value = -1500623

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 552.7
$ws.Range("I113").Value = 365.57144
$ws.Range("J113").Value = 653.46155
$ws.Range("K113").Value = 1096.71432
$ws.Range("L113").Value = 1960.38465
$ws.Range("M113").Value = 1073.28568
$ws.Range("N113").Value = -6300.38465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 4528
$ws.Range("I114").Value = 4528
$ws.Range("K114").Value = 13584
$ws.Range("M114").Value = -10330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4204.7
$ws.Range("I140").Value = 4204.7
$ws.Range("K140").Value = 12614.1
$ws.Range("M140").Value = -7434.099999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1598.8889
$ws.Range("J22").Value = 1550
$ws.Range("L22").Value = 1550
$ws.Range("N22").Value = -2140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1598.8889
$ws.Range("J27").Value = 1550
$ws.Range("L27").Value = 1550
$ws.Range("N27").Value = -1764

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 10000
$ws.Range("I34").Value = 10000
$ws.Range("K34").Value = 10000
$ws.Range("M34").Value = -9828

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4602.84
$ws.Range("I40").Value = 3503.4
$ws.Range("K40").Value = 3503.4
$ws.Range("M40").Value = -3367.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4367.3
$ws.Range("I61").Value = 4379.8237
$ws.Range("K61").Value = 4379.8237
$ws.Range("M61").Value = -4177.8237

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4367.3
$ws.Range("I113").Value = 4379.8237
$ws.Range("K113").Value = 4379.8237
$ws.Range("M113").Value = -2209.8237

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5500.875
$ws.Range("I132").Value = 4801.4
$ws.Range("K132").Value = 14404.2
$ws.Range("M132").Value = -11874.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1941.8334
$ws.Range("J96").Value = 2080.875
$ws.Range("L96").Value = 2080.875
$ws.Range("N96").Value = -4826.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 333355000
$ws.Range("J125").Value = 333355000
$ws.Range("L125").Value = 333355000
$ws.Range("N125").Value = -333364840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5152.8
$ws.Range("J126").Value = 5642.4287
$ws.Range("L126").Value = 16927.2861
$ws.Range("N126").Value = -21867.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2218.2068
$ws.Range("I132").Value = 2160.2964
$ws.Range("K132").Value = 6480.889200000001
$ws.Range("M132").Value = -3950.889200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2554.5356
$ws.Range("I136").Value = 2396.4375
$ws.Range("K136").Value = 7189.3125
$ws.Range("M136").Value = -4639.3125
